# RPA datasets push 2023-12-28
# Target sheet: "02_38커뮤니케이션(최근일자기준)" (the IPO demand-forecast table, columns A:F, rows 1:21)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("02_38커뮤니케이션(최근일자기준)")

# Insert a new data row right after the header+first row (new row 3), shifting
# the existing rows 3:21 down to 4:22.
$ws.Range("A3").EntireRow.Insert()

# Populate the newly inserted row with the new IPO entry (IBKS스팩24호).
$ws.Range("A3").Value = "IBKS스팩24호"
$ws.Range("B3").Value = "2024.01.17~01.18"
$ws.Range("C3").Value = "2,000~2,000"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = 2147483647
$ws.Range("F3").Value = "아이비케이투자증권"

# Drop the oldest entry (was row 21, now shifted to row 22 - 와이바이오로직스)
# so the table keeps its original 20-data-row size.
$ws.Range("A22").EntireRow.Delete()
